$d = $word.ActiveDocument

# Update the payment-schedule paragraph: the payment window changes from
# "Monday to Saturday, 8:30 AM - 5:30 PM" to "Monday to Friday 9:00 AM -
# 5:00 PM and Saturday 9:00 AM - 2:00 PM", and the appointment notice period
# gains "al menos" (at least).
$oldText = "LOS PAGOS DEBERÁN REALIZARSE DE LUNES A SÁBADO, ENTRE LAS 8:30 A.M. Y LAS 5:30 P.M. PARA EFECTUAR UN PAGO EN DOMINGO, SERÁ INDISPENSABLE PROGRAMAR UNA CITA CON TRES DÍAS DE ANTICIPACIÓN. CADA PAGO DEBERÁ SER NOTIFICADO Y CONFIRMADO AL NÚMERO TELEFÓNICO 951 189 9298."
$newText = "LOS PAGOS DEBERÁN REALIZARSE DE LUNES A VIERNES, EN UN HORARIO DE 9:00 A.M. A 5:00 P.M., Y EN SÁBADO DE 9:00 A. M. A 2:00 P. M. PARA EFECTUAR UN PAGO EN DOMINGO, SERÁ INDISPENSABLE PROGRAMAR UNA CITA CON AL MENOS TRES DÍAS DE ANTICIPACIÓN. CADA PAGO DEBERÁ SER NOTIFICADO Y CONFIRMADO AL NÚMERO TELEFÓNICO 951 189 9298."

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Output "ReplacedPaymentScheduleText=$found"
